$wb = $excel.ActiveWorkbook

# "User" sheet: update the course text and move the active selection
$wsUser = $wb.Worksheets.Item("User")
$wsUser.Range("D2").Value = "Web Development"

# Make "User" the active sheet/tab and move its selection to F14
$wsUser.Activate()
$wsUser.Range("F14").Select()
